$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source row for 2007 年 (row 2) was removed from the sheet and all
# subsequent rows (2010, 2012, 2015, 2017) shifted up by one.
$ws.Rows.Item(2).Delete()
